# Refresh cryptocurrency price/volume snapshot (coinranking.com feed).
# Price/volume cells are stored as plain text (not numbers) in this sheet,
# so numeric-looking values that Excel would otherwise auto-convert are
# entered with a leading apostrophe to force a text literal, matching the
# original inline-string cell type.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "97.603.90"
$ws.Range("E2").Value = "  +0.84%  "
$ws.Range("D3").Value = "3.295.39"
$ws.Range("E3").Value = "  -0.94%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.20%  "
$ws.Range("D5").Value = "'254.31"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("D6").Value = "'621.46"
$ws.Range("E6").Value = "  +0.04%  "
$ws.Range("D7").Value = "'1.44"
$ws.Range("E7").Value = "  +28.31%  "
$ws.Range("D8").Value = "'0.398"
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "'0.999"
$ws.Range("E9").Value = "  -0.18%  "
$ws.Range("D10").Value = "'0.904"
$ws.Range("E10").Value = "  +14.89%  "
$ws.Range("D11").Value = "3.294.14"
$ws.Range("E11").Value = "  -0.91%  "
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").Value = "'38.60"
$ws.Range("E13").Value = "  +9.91%  "
$ws.Range("D14").Value = "97.320.46"
$ws.Range("E14").Value = "  +0.40%  "
$ws.Range("D15").Value = "'0.0000246"
$ws.Range("E15").Value = "  -0.36%  "
$ws.Range("D16").Value = "3.913.69"
$ws.Range("E16").Value = "  -0.98%  "
$ws.Range("E17").Value = "  -0.78%  "
$ws.Range("D18").Value = "3.299.40"
$ws.Range("E18").Value = "  -1.34%  "
$ws.Range("D19").Value = "'3.49"
$ws.Range("E19").Value = "  -2.70%  "
$ws.Range("D20").Value = "'15.07"
$ws.Range("E20").Value = "  +0.22%  "
$ws.Range("D21").Value = "'6.16"
$ws.Range("E21").Value = "  +5.86%  "
$ws.Range("D22").Value = "'477.07"
$ws.Range("E22").Value = "  -1.44%  "
$ws.Range("D23").Value = "'9.38"
$ws.Range("E23").Value = "  +0.69%  "
$ws.Range("D24").Value = "'0.0000202"
$ws.Range("E24").Value = "  -2.87%  "
$ws.Range("D25").Value = "'5.57"
$ws.Range("E25").Value = "  -1.42%  "
$ws.Range("D26").Value = "'87.94"
$ws.Range("E26").Value = "  -0.18%  "
$ws.Range("D27").Value = "'11.76"
$ws.Range("E27").Value = "  -2.72%  "
$ws.Range("D28").Value = "'0.298"
$ws.Range("E28").Value = "  +24.26%  "
$ws.Range("D29").Value = "3.467.01"
$ws.Range("E29").Value = "  -1.99%  "
$ws.Range("D30").Value = "'0.998"
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("D31").Value = "'0.185"
$ws.Range("E31").Value = "  +1.78%  "
$ws.Range("E32").Value = "  +8.50%  "
$ws.Range("D33").Value = "'9.83"
$ws.Range("E33").Value = "  +6.19%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "'27.50"
$ws.Range("E35").Value = "  +0.89%  "
$ws.Range("E36").Value = "  -3.59%  "
$ws.Range("D37").Value = "'7.14"
$ws.Range("E37").Value = "  -4.11%  "
$ws.Range("E38").Value = "  +0.17%  "
$ws.Range("D39").Value = "'24.78"
$ws.Range("E39").Value = "  +0.51%  "
$ws.Range("B40").Value = "PolygonEcosystemToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D40").Value = "'0.453"
$ws.Range("E40").Value = "  +0.90%  "
$ws.Range("B41").Value = "Bittensor"
$ws.Range("C41").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D41").Value = "'487.61"
$ws.Range("E41").Value = "  -2.52%  "
$ws.Range("D42").Value = "'3.69"
$ws.Range("E42").Value = "  +6.09%  "
$ws.Range("E43").Value = "  -2.91%  "
$ws.Range("D44").Value = "'0.798"
$ws.Range("E44").Value = "  +1.59%  "
$ws.Range("E45").Value = "  -0.02%  "
$ws.Range("E46").Value = "  -5.30%  "
$ws.Range("D47").Value = "'158.18"
$ws.Range("E47").Value = "  -1.77%  "
$ws.Range("E48").Value = "  -1.40%  "
$ws.Range("D49").Value = "'0.839"
$ws.Range("E49").Value = "  +5.88%  "
$ws.Range("B50").Value = "Optimism"
$ws.Range("C50").Value = "https://coinranking.com/coin/n1p-s_gm1+optimism-op"
$ws.Range("D50").Value = "'2.20"
$ws.Range("E50").Value = "  +5.71%  "
$ws.Range("E51").Value = "  +1.81%  "
